# Applies the cryptos-list price/volume refresh described by the commit
# ("Updated cryptos list on Sat Sep 23 23:27:24 UTC 2023 with GitHub Actions").
# Column D = Price, Column E = Volume(1h); every value in the sheet is stored
# as literal text (not a Number/Percentage cell), so each update below writes
# the exact replacement text for one cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, [string]$cellRef, [string]$text, [bool]$looksNumeric) {
    $range = $ws.Range($cellRef)
    if ($looksNumeric) {
        # Plain '.Value = "211.59"' would make Excel coerce the cell to a
        # Number (risking float drift, e.g. 19.53 -> 19.530000000000001, and
        # losing formatting like trailing zeros/dot-separators). Force text
        # entry the same way Excel's UI does it: a leading apostrophe (quote
        # prefix), then put the style back to Normal so the cell ends up with
        # no formatting override - matching the original, unstyled text cell.
        $range.Value = "'" + $text
        $range.Style = 'Normal'
    } else {
        $range.Value = $text
    }
}

Set-TextValue $ws 'D2' '26.714.08' $false
Set-TextValue $ws 'E2' '  +0.24%  ' $false
Set-TextValue $ws 'D3' '1.599.92' $false
Set-TextValue $ws 'E3' '  +0.20%  ' $false
Set-TextValue $ws 'E4' '  +0.36%  ' $false
Set-TextValue $ws 'D5' '211.59' $true
Set-TextValue $ws 'E5' '  +0.01%  ' $false
Set-TextValue $ws 'E6' '  -0.69%  ' $false
Set-TextValue $ws 'E7' '  +0.36%  ' $false
Set-TextValue $ws 'E8' '  +0.11%  ' $false
Set-TextValue $ws 'D9' '0.247' $true
Set-TextValue $ws 'E9' '  +0.76%  ' $false
Set-TextValue $ws 'D10' '19.53' $true
Set-TextValue $ws 'E10' '  +0.02%  ' $false
Set-TextValue $ws 'E11' '  +0.53%  ' $false
Set-TextValue $ws 'D12' '1.824.98' $false
Set-TextValue $ws 'E12' '  +0.22%  ' $false
Set-TextValue $ws 'D13' '1.609.09' $false
Set-TextValue $ws 'E13' '  +4.04%  ' $false
Set-TextValue $ws 'E14' '  +0.57%  ' $false
Set-TextValue $ws 'E15' '  +0.22%  ' $false
Set-TextValue $ws 'D16' '65.33' $true
Set-TextValue $ws 'E16' '  +1.41%  ' $false
Set-TextValue $ws 'D17' '26.687.30' $false
Set-TextValue $ws 'E17' '  +0.24%  ' $false
Set-TextValue $ws 'D18' '0.0₃0754' $false
Set-TextValue $ws 'E18' '  +2.98%  ' $false
Set-TextValue $ws 'D19' '7.22' $true
Set-TextValue $ws 'E19' '  +3.84%  ' $false
Set-TextValue $ws 'E20' '  +0.35%  ' $false
Set-TextValue $ws 'D21' '209.30' $true
Set-TextValue $ws 'E21' '  +0.31%  ' $false
Set-TextValue $ws 'E22' '  +0.50%  ' $false
Set-TextValue $ws 'E24' '  +0.67%  ' $false
Set-TextValue $ws 'D25' '142.27' $true
Set-TextValue $ws 'E25' '  -1.97%  ' $false
Set-TextValue $ws 'E26' '  +0.34%  ' $false
Set-TextValue $ws 'E27' '  -0.59%  ' $false
Set-TextValue $ws 'E28' '  -0.12%  ' $false
Set-TextValue $ws 'E29' '  +0.53%  ' $false
Set-TextValue $ws 'D30' '0.0522' $true
Set-TextValue $ws 'E30' '  +3.08%  ' $false
Set-TextValue $ws 'E31' '  -0.27%  ' $false
Set-TextValue $ws 'E32' '  +0.64%  ' $false
Set-TextValue $ws 'E33' '  +1.87%  ' $false
Set-TextValue $ws 'D34' '1.292.19' $false
Set-TextValue $ws 'E34' '  +1.04%  ' $false
Set-TextValue $ws 'D35' '0.624' $true
Set-TextValue $ws 'E35' '  -5.08%  ' $false
Set-TextValue $ws 'E36' '  +0.98%  ' $false
Set-TextValue $ws 'E37' '  +0.26%  ' $false
Set-TextValue $ws 'E38' '  +0.04%  ' $false
Set-TextValue $ws 'D39' '1.11' $true
Set-TextValue $ws 'E39' '  +20.12%  ' $false
Set-TextValue $ws 'D40' '0.825' $true
Set-TextValue $ws 'E40' '  -2.31%  ' $false
Set-TextValue $ws 'D41' '5.42' $true
Set-TextValue $ws 'E41' '  -0.88%  ' $false
Set-TextValue $ws 'D42' '2.20' $true
Set-TextValue $ws 'E42' '  -0.20%  ' $false
Set-TextValue $ws 'D43' '0.783' $true
Set-TextValue $ws 'E43' '  -0.40%  ' $false
Set-TextValue $ws 'D44' '63.18' $true
Set-TextValue $ws 'E44' '  -2.13%  ' $false
Set-TextValue $ws 'D45' '1.736.53' $false
Set-TextValue $ws 'D46' '91.32' $true
Set-TextValue $ws 'E46' '  +1.54%  ' $false
Set-TextValue $ws 'E47' '  -1.63%  ' $false
Set-TextValue $ws 'E48' '  +0.73%  ' $false
Set-TextValue $ws 'E49' '  -1.22%  ' $false
Set-TextValue $ws 'E50' '  +0.56%  ' $false
Set-TextValue $ws 'E51' '  +0.38%  ' $false
